# Update "CLAN CAPITAL" sheet - add new weekly event column Q (20/02/2026),
# add new member rows, and update a few member names/values - 2026-02-21

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CLAN CAPITAL")

# --- New header for column Q (was a placeholder "Colonna15") ---
$ws.Range("Q1").Value = "20/02/2026"

# --- Fill in column Q results for all existing members (rows 2-43) ---
$qResults = @(
    @(2, 0), @(3, 0), @(4, 6), @(5, 6), @(6, 0), @(7, 0), @(8, 0), @(9, 0),
    @(10, 4), @(11, 0), @(12, 0), @(13, 2), @(14, 0), @(15, 6), @(16, 6),
    @(17, 0), @(18, 0), @(19, 0), @(20, 0), @(21, 6), @(22, 0), @(23, 3),
    @(24, 0), @(25, 0), @(26, 0), @(27, 0), @(28, 6), @(29, 6), @(30, 4),
    @(31, 6), @(32, 0), @(33, 6), @(34, 6), @(35, 0), @(36, 0), @(37, 0),
    @(38, 6), @(39, 0), @(40, 6), @(41, 0), @(42, 6), @(43, 0)
)
foreach ($pair in $qResults) {
    $r = $pair[0]
    $v = $pair[1]
    $ws.Cells.Item($r, 17).Value = $v
    $ws.Cells.Item($r, 17).Style = "Normal"
}

# --- Row 44: name shifts up the roster (Tom -> fede61mito) ---
$ws.Range("A44").Value = "fede61mito"
$ws.Range("Q44").Value = 0
$ws.Range("Q44").Style = "Normal"

# --- Row 45: name shifts up (fede61mito -> Dasters79), and its score updates ---
$ws.Range("A45").Value = "Dasters79"
$ws.Range("P45").Value = 1
$ws.Range("Q45").Value = 2
$ws.Range("Q45").Style = "Normal"

# --- Row 46: name shifts up (Dasters79 -> Amir✴), new join date, old P score removed ---
$ws.Range("A46").Value = "Amir✴"
$ws.Range("B46").Value = "21/02/2026"
$ws.Range("P46").ClearContents()
$ws.Range("Q46").Value = 0
$ws.Range("Q46").Style = "Normal"

# --- New members joining 21/02/2026 (rows 47-51) ---
$newMembers = @(
    @(47, "Artur"),
    @(48, "Xx_Herman_xX"),
    @(49, "dibba10"),
    @(50, "Anto"),
    @(51, "Michele")
)
foreach ($member in $newMembers) {
    $r = $member[0]
    $name = $member[1]

    $ws.Range("A${r}").Value = $name
    $ws.Range("B${r}").Value = "21/02/2026"
    $ws.Range("C${r}").Formula = "=ROUND(AVERAGE(D${r}:AH${r}), 0)"

    # Bring over the same blank-but-bordered K:N formatting used by the rest of the roster
    $ws.Range("K46:N46").Copy()
    $ws.Range("K${r}:N${r}").PasteSpecial(-4122)

    $ws.Range("Q${r}").Value = 0
    $ws.Range("Q${r}").Style = "Normal"
}
